$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.646.20'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").Value = '1.593.08'
$ws.Range("E3").Value = '  -1.75%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.24'
$ws.Range("E5").Value = '  -1.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.512'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.245'
$ws.Range("E9").Value = '  -1.88%  '
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0836'
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '1.815.59'
$ws.Range("E12").Value = '  -1.80%  '
$ws.Range("D13").Value = '1.592.44'
$ws.Range("E13").Value = '  -1.56%  '
$ws.Range("E14").Value = '  -2.47%  '
$ws.Range("E15").Value = '  -3.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.67'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").Value = '26.618.07'
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '209.22'
$ws.Range("E19").Value = '  -2.63%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("E22").Value = '  -2.38%  '
$ws.Range("E23").Value = '  -2.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.90'
$ws.Range("E24").Value = '  -1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.52'
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.12'
$ws.Range("E27").Value = '  -4.11%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.28'
$ws.Range("E29").Value = '  -1.47%  '
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.15'
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("E32").Value = '  -2.90%  '
$ws.Range("E33").Value = '  -4.78%  '
$ws.Range("E34").Value = '  -2.84%  '
$ws.Range("D35").Value = '1.291.95'
$ws.Range("E35").Value = '  -3.22%  '
$ws.Range("E37").Value = '  -5.40%  '
$ws.Range("E38").Value = '  -2.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.838'
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.790'
$ws.Range("E41").Value = '  -0.49%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.35'
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.19'
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.44'
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").Value = '1.728.55'
$ws.Range("E45").Value = '  -1.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.895'
$ws.Range("E46").Value = '  +4.84%  '
$ws.Range("E47").Value = '  -0.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.64'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0981'
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0503'
$ws.Range("E50").Value = '  -1.63%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.27%  '
